# Atualizado por script em 08-11-2023 14:45
#
# The source scraper re-ran and, for several match days, produced the two
# fixtures of that round in a different order than before (the row's
# metadata in columns A-E - index/country/competition/season/date - stayed
# put, but the match details in columns F-V got shuffled between the two
# rows). On top of that a brand-new fixture (Persis Solo x PSS Sleman) was
# appended as row 164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row-Details($rowA, $rowB) {
    $rangeA = $ws.Range("F$rowA`:V$rowA")
    $rangeB = $ws.Range("F$rowB`:V$rowB")
    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()
    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Pairs of rows whose match details (columns F..V) were swapped.
Swap-Row-Details 19 20
Swap-Row-Details 37 38
Swap-Row-Details 49 50
Swap-Row-Details 51 52
Swap-Row-Details 84 85
Swap-Row-Details 139 140
Swap-Row-Details 162 163

# New fixture appended at the end of the sheet.
$newRow = 164

# Copy formatting of the columns that carry special styling (index / date)
# from the previous row so the new row matches the sheet's layout.
$ws.Range("A163").Copy() | Out-Null
$ws.Range("A$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E163").Copy() | Out-Null
$ws.Range("E$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A$newRow").Value = 163
$ws.Range("B$newRow").Value = "indonesia"
$ws.Range("C$newRow").Value = "liga-1"
$ws.Range("D$newRow").Value = "2023-2024"
$ws.Range("E$newRow").Value = 45238.54166666666
$ws.Range("F$newRow").Value = "Persis Solo"
$ws.Range("G$newRow").Value = 1
$ws.Range("H$newRow").Value = "PSS Sleman"
$ws.Range("I$newRow").Value = 1
$ws.Range("J$newRow").Value = 1.57
$ws.Range("K$newRow").Value = "07/11/2023 01:12"
$ws.Range("L$newRow").Value = 1.68
$ws.Range("M$newRow").Value = "08/11/2023 12:58"
$ws.Range("N$newRow").Value = 3.99
$ws.Range("O$newRow").Value = "07/11/2023 01:12"
$ws.Range("P$newRow").Value = 3.91
$ws.Range("Q$newRow").Value = "08/11/2023 12:59"
$ws.Range("R$newRow").Value = 4.68
$ws.Range("S$newRow").Value = "07/11/2023 01:12"
$ws.Range("T$newRow").Value = 4.88
$ws.Range("U$newRow").Value = "08/11/2023 12:59"
$ws.Range("V$newRow").Value = "https://www.betexplorer.com/football/indonesia/liga-1/persis-solo-pss-sleman/GYmhOpdI/"

$wb.Save()
